$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9408089518547058
$ws.Range("B1").Value = 1.975756764411926
$ws.Range("C1").Value = 7.504281520843506
$ws.Range("D1").Value = 2.724705457687378
$ws.Range("E1").Value = 1.366336345672607
